# fix line endings for bash script examples
#
# Several multi-line "bash" verbatim examples use a trailing backslash
# ("\") on every line except the last to indicate the command continues
# on the next line. The "--remote"/"--remotes ..." line was missing its
# trailing backslash in a few of these examples (wherever it is not the
# final line of the script), so add it back — but only where the script
# actually continues afterwards.

$d = $word.ActiveDocument

$oldRemote  = " --remote https://linker.bio,https://zenodo.org"
$newRemote  = " --remote https://linker.bio,https://zenodo.org\"
$oldRemotes = " --remotes https://linker.bio,https://zenodo.org"
$newRemotes = " --remotes https://linker.bio,https://zenodo.org\"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $trimmed = $text.TrimEnd([char]13, [char]7, [char]10)

    # If the URL is the very last thing in the paragraph, this is the
    # final line of its script block, so it must stay unchanged.
    if ($trimmed.EndsWith($oldRemote) -or $trimmed.EndsWith($oldRemotes)) {
        continue
    }

    $rng = $para.Range
    $rng.Find.Execute($oldRemote, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $newRemote, 2) | Out-Null

    $rng2 = $para.Range
    $rng2.Find.Execute($oldRemotes, $true, $false, $false, $false, $false, `
                        $true, 1, $false, $newRemotes, 2) | Out-Null
}
